# The diff moves the "Meta description: ..." paragraph that originally sat
# right under the H1 title down to the very end of the document (just before
# the final "Create a feature image..." paragraph), splitting it into two
# separate paragraphs there:
#   - a bold "Play Cash Falls Island Bounty for Free - Review" paragraph
#   - an italic paragraph with the old description text (minus the
#     "Meta description: " prefix), which *replaces* the old
#     "Create a feature image..." paragraph text.
# The original "Meta description" paragraph at the top is then removed.

$d = $word.ActiveDocument

# --- Step 1: duplicate the "Meta description" paragraph (paragraph 2) to
#     just before the final paragraph. ---

$metaPara = $d.Paragraphs.Item(2)
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)

$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.FormattedText = $metaPara.Range.FormattedText

# The pasted paragraph is now paragraph number $n (the old last paragraph got
# pushed one slot further down, to $n + 1).
$newPara = $d.Paragraphs.Item($n)
$newRange = $newPara.Range

# Turn "Meta description" (the bold run) into the new bold heading text.
$renamed = $newRange.Find.Execute("Meta description", $true, $false, $false, `
    $false, $false, $true, 1, $false, `
    "Play Cash Falls Island Bounty for Free - Review", 2)
if (-not $renamed) {
    throw "Could not find 'Meta description' run to rename"
}

# Drop the old ": Read our neutral review ..." text that used to trail it,
# leaving just the new bold paragraph behind.
$newPara = $d.Paragraphs.Item($n)
$descRange = $newPara.Range.Duplicate
$found = $descRange.Find.Execute(": Read our neutral review of Cash Falls Island Bounty. Play this online slot game for free and find out about its gameplay features, symbols, and much more.")
if ($found) {
    $descRange.Delete()
}

# --- Step 2: replace the final paragraph's (italic) text with the
#     description, now that it has moved here. ---

$n2 = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($n2)
$finalRange = $finalPara.Range
$replaced = $finalRange.Find.Execute("Create a feature image fitting for the game ""Cash Falls Island Bounty"". Please design a cartoon-style image featuring a happy Maya warrior with glasses. The Maya warrior should be depicted on a Caribbean beach, with the ocean in the background and a treasure chest overflowing with gold coins and jewels in the foreground. The warrior should be holding a large gold coin with the game's logo engraved on it, and have a big smile on his face. The cartoon-style image should be bright, colorful, and inviting.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Read our neutral review of Cash Falls Island Bounty. Play this online slot game for free and find out about its gameplay features, symbols, and much more.", 2)
if (-not $replaced) {
    throw "Could not find the old feature-image-prompt paragraph text"
}

# --- Step 3: remove the original "Meta description" paragraph from the top
#     of the document (it now lives at the end instead). ---

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$metaRange.Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
